$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3791.1667
$ws.Range("J18").Value = 5999.6665
$ws.Range("L18").Value = 5999.6665
$ws.Range("N18").Value = -6567.6665

$ws.Range("H86").Value = 12503597
$ws.Range("J86").Value = 16671059
$ws.Range("L86").Value = 16671059
$ws.Range("N86").Value = -16673305

$ws.Range("H88").Value = 7219.25
$ws.Range("J88").Value = 9683.723
$ws.Range("L88").Value = 9683.723
$ws.Range("N88").Value = -10495.723

$ws.Range("H89").Value = 12503597
$ws.Range("J89").Value = 16671059
$ws.Range("L89").Value = 83355295
$ws.Range("N89").Value = -83366527

$ws.Range("H91").Value = 7219.25
$ws.Range("J91").Value = 9683.723
$ws.Range("L91").Value = 9683.723
$ws.Range("N91").Value = -12491.723

$ws.Range("H100").Value = 4867.4287
$ws.Range("I100").Value = 2575.3333
$ws.Range("J100").Value = 6586.5
$ws.Range("K100").Value = 2575.3333
$ws.Range("L100").Value = 6586.5
$ws.Range("M100").Value = -2034.3333
$ws.Range("N100").Value = -7668.5

$ws.Range("H116").Value = 243116.4
$ws.Range("I116").Value = 10914.833
$ws.Range("J116").Value = 397917.44
$ws.Range("K116").Value = 10914.833
$ws.Range("L116").Value = 397917.44
$ws.Range("M116").Value = -7472.833000000001
$ws.Range("N116").Value = -404801.44

$ws.Range("H125").Value = 20000828
$ws.Range("I125").Value = 1035
$ws.Range("J125").Value = 100000000
$ws.Range("K125").Value = 9315
$ws.Range("L125").Value = 900000000
$ws.Range("M125").Value = -6855
$ws.Range("N125").Value = -900004920

$ws.Range("H131").Value = 3477.0588
$ws.Range("I131").Value = 1535.2727
$ws.Range("K131").Value = 4605.8181
$ws.Range("M131").Value = 434.1818999999996

$ws.Range("H132").Value = 3520620.5
$ws.Range("I132").Value = 3974636.2
$ws.Range("K132").Value = 11923908.6
$ws.Range("M132").Value = -11921378.6

$ws.Range("H138").Value = 4254.9067
$ws.Range("I138").Value = 2712.8
$ws.Range("J138").Value = 4457.816
$ws.Range("K138").Value = 8138.400000000001
$ws.Range("L138").Value = 13373.448
$ws.Range("M138").Value = -2998.400000000001
$ws.Range("N138").Value = -23653.448

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 3781.5
$ws.Range("I17").Value = 2008
$ws.Range("J17").Value = 5555
$ws.Range("K17").Value = 2008
$ws.Range("L17").Value = 5555
$ws.Range("M17").Value = -1835
$ws.Range("N17").Value = -5901

$ws.Range("H32").Value = 2752.8357
$ws.Range("I32").Value = 2776.3845
$ws.Range("K32").Value = 2776.3845
$ws.Range("M32").Value = -2489.3845

$ws.Range("H61").Value = 906302.75
$ws.Range("I61").Value = 1015309.2
$ws.Range("K61").Value = 1015309.2
$ws.Range("M61").Value = -1015097.2

$ws.Range("H122").Value = 2267.6223
$ws.Range("I122").Value = 1480.1515
$ws.Range("J122").Value = 4433.1665
$ws.Range("K122").Value = 4440.4545
$ws.Range("L122").Value = 13299.4995
$ws.Range("M122").Value = -1990.4545
$ws.Range("N122").Value = -18199.4995

$ws.Range("H132").Value = 427680.4
$ws.Range("I132").Value = 489866.56
$ws.Range("J132").Value = 2741.6667
$ws.Range("K132").Value = 1469599.68
$ws.Range("L132").Value = 8225.000100000001
$ws.Range("M132").Value = -1467069.68
$ws.Range("N132").Value = -13285.0001

$ws.Range("H136").Value = 906302.75
$ws.Range("I136").Value = 1015309.2
$ws.Range("K136").Value = 3045927.6
$ws.Range("M136").Value = -3043377.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 23723.555
$ws.Range("I26").Value = 23723.555
$ws.Range("K26").Value = 23723.555
$ws.Range("M26").Value = -23431.555

$ws.Range("H96").Value = 18140.125
$ws.Range("I96").Value = 18140.125
$ws.Range("K96").Value = 18140.125
$ws.Range("M96").Value = -15394.125

$ws.Range("H105").Value = 3305.25
$ws.Range("I105").Value = 2914.3845
$ws.Range("K105").Value = 2914.3845
$ws.Range("M105").Value = -1167.3845

$ws.Range("H108").Value = 99622.664
$ws.Range("J108").Value = 99622.664
$ws.Range("L108").Value = 99622.664
$ws.Range("N108").Value = -107302.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H127").Value = 29000
$ws.Range("J127").Value = 29000
$ws.Range("L127").Value = 29000
$ws.Range("N127").Value = -38920

$ws.Range("H132").Value = 466665.4
$ws.Range("I132").Value = 494148.84
$ws.Range("K132").Value = 1482446.52
$ws.Range("M132").Value = -1479916.52

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 3029.3157
$ws.Range("J140").Value = 4118.9165
$ws.Range("L140").Value = 12356.7495
$ws.Range("N140").Value = -22716.7495

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 70261
$ws.Range("J39").Value = 70261
$ws.Range("L39").Value = 70261
$ws.Range("N39").Value = -71325

$ws.Range("H122").Value = 4329.515
$ws.Range("I122").Value = 5308.3335
$ws.Range("K122").Value = 15925.0005
$ws.Range("M122").Value = -13475.0005

$ws.Range("H132").Value = 318530.9
$ws.Range("I132").Value = 348905.53
$ws.Range("K132").Value = 1046716.59
$ws.Range("M132").Value = -1044186.59

$ws.Range("H141").Value = 37085.6
$ws.Range("J141").Value = 37085.6
$ws.Range("L141").Value = 37085.6
$ws.Range("N141").Value = -47445.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 18000
$ws.Range("J3").Value = 18000
$ws.Range("L3").Value = 18000
$ws.Range("N3").Value = -18224

$ws.Range("H15").Value = 18000
$ws.Range("J15").Value = 18000
$ws.Range("L15").Value = 18000
$ws.Range("N15").Value = -18340

$ws.Range("H100").Value = 14540.375
$ws.Range("I100").Value = 2391.5
$ws.Range("K100").Value = 2391.5
$ws.Range("M100").Value = -1850.5

$ws.Range("H136").Value = 4148.1904
$ws.Range("I136").Value = 3450.7222
$ws.Range("K136").Value = 10352.1666
$ws.Range("M136").Value = -7802.1666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 250017500
$ws.Range("I4").Value = 19995
$ws.Range("J4").Value = 333349980
$ws.Range("K4").Value = 19995
$ws.Range("L4").Value = 333349980
$ws.Range("M4").Value = -19882
$ws.Range("N4").Value = -333350206

$ws.Range("H107").Value = 7498.8335
$ws.Range("I107").Value = 10001.5
$ws.Range("J107").Value = 6247.5
$ws.Range("K107").Value = 30004.5
$ws.Range("L107").Value = 18742.5
$ws.Range("M107").Value = -28084.5
$ws.Range("N107").Value = -22582.5

$ws.Range("H126").Value = 3115.5386
$ws.Range("I126").Value = 2736.818
$ws.Range("K126").Value = 8210.454000000002
$ws.Range("M126").Value = -5740.454000000002
